# Generate Report for Archive
# - Status text changes from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns, E & F) and on
#   the per-locale sheets ("zh-cn", "de-de") in their "Status" column (C).
# - The Status columns shrink to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($cellRef in @("E2", "F2", "E3", "F3")) {
    $wsOverview.Range($cellRef).Value = $newStatus
}
# Shrink the now-narrower Status columns to their best-fit width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766

# --- Per-locale sheets: column C ("Status"), rows 2-3 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in @("C2", "C3")) {
        $ws.Range($cellRef).Value = $newStatus
    }
    $ws.Columns.Item(3).ColumnWidth = 12.576851254417766
}
